# Update the cached text of the "datetimeFigureOut" date placeholder field
# from 4/18/17 to 4/20/17 across the slide master, every slide layout, and
# the notes master.

$p = $ppt.ActivePresentation
$newDate = "4/20/17"
$oldDate = "4/18/17"

function Update-DatePlaceholderShapes($shapes, $newText, $oldText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.Name.StartsWith("Date Placeholder")) {
            if ($sh.TextFrame.TextRange.Text -eq $oldText) {
                $sh.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

# Slide master date placeholder
Update-DatePlaceholderShapes $p.SlideMaster.Shapes $newDate $oldDate

# Every slide layout that hangs off the (single) slide master
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Update-DatePlaceholderShapes $layout.Shapes $newDate $oldDate
}

# Notes master date placeholder -- updated through the HeadersFooters /
# DateAndTime object, since the NotesMaster.Shapes collection does not
# route text writes to the correct underlying shape.
$nmDateAndTime = $p.NotesMaster.HeadersFooters.DateAndTime
$nmDateAndTime.Text = $newDate
